# add libraries poi-3.17 xmlbeans-2.6.0 variable declaration of input data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows 12-16 (A=Date, B=vechical name, C=S.km, D=C.km, E=Rent,
# F=D.rate, G=d.quantity, H=Milage, I=M.rate, J=Total, K=M.Des)
$rows = @(
    @{ Row=12; A="14-03-2018"; B="L3";  C=55000.0;  D=60000.0;  E=6000.0;  F=66.0;  G=150.0;  H=33.0;  I=0.0;    J=19600.0;    K="0" },
    @{ Row=13; A="06-03-2018"; B="L2";  C=654.0;    D=654.0;    E=654.0;   F=654.0; G=654.0;  H=0.0;   I=6546.0; J=-46.0;      K="46" },
    @{ Row=14; A="08-03-2018"; B="L1";  C=64.0;     D=654.0;    E=5465.0;  F=654.0; G=6546.0; H=0.0;   I=465.0;  J=-707028.0;  K="465" },
    @{ Row=15; A="06-03-2018"; B="L1";  C=4684.0;   D=6868.0;   E=468.0;   F=84.0;  G=468.0;  H=4.0;   I=468.0;  J=229266.0;   K="48" },
    @{ Row=16; A="05-03-2018"; B="L1";  C=64.0;     D=65465.0;  E=465.0;   F=654.0; G=465.0;  H=140.0; I=4654.0; J=-264309.0;  K="65" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    # Date (A), vehicle-name (B) and M.Des (K) are all plain text in the
    # source data, including numeric- and date-looking values ("0", "46",
    # "06-03-2018", ...). Force text storage with a leading apostrophe so
    # Excel doesn't auto-convert them to numbers/dates, then clear the
    # resulting quote-prefix style so the cells stay unstyled like the
    # rest of the data rows.
    $ws.Cells.Item($rowNum, 1).Value = "'" + $r.A
    $ws.Cells.Item($rowNum, 1).Style = "Normal"

    $ws.Cells.Item($rowNum, 2).Value = "'" + $r.B
    $ws.Cells.Item($rowNum, 2).Style = "Normal"

    $ws.Cells.Item($rowNum, 3).Value = $r.C
    $ws.Cells.Item($rowNum, 4).Value = $r.D
    $ws.Cells.Item($rowNum, 5).Value = $r.E
    $ws.Cells.Item($rowNum, 6).Value = $r.F
    $ws.Cells.Item($rowNum, 7).Value = $r.G
    $ws.Cells.Item($rowNum, 8).Value = $r.H
    $ws.Cells.Item($rowNum, 9).Value = $r.I
    $ws.Cells.Item($rowNum, 10).Value = $r.J

    $ws.Cells.Item($rowNum, 11).Value = "'" + $r.K
    $ws.Cells.Item($rowNum, 11).Style = "Normal"
}
